$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = -0.1652581622154127
$ws.Range("D2").Value = -0.002523758263964982
$ws.Range("E2").Value = 0.0002931071949994339
$ws.Range("F2").Value = 0.001124234086298129
$ws.Range("G2").Value = -0.00127753765741268
$ws.Range("H2").Value = -0.0007990859124778034
$ws.Range("I2").Value = 0.001177632623395438
$ws.Range("J2").Value = 0.0004195049823813824
$ws.Range("K2").Value = -0.00219475635291615
$ws.Range("L2").Value = -0.001166882514834089
$ws.Range("M2").Value = 0.001634561894822562
$ws.Range("N2").Value = 0.0008368912573226979
$ws.Range("O2").Value = -0.001468032878502348
$ws.Range("P2").Value = 0.001248231948587792
$ws.Range("Q2").Value = -0.0001451763067674667
$ws.Range("R2").Value = 0.002522203803982969
$ws.Range("B3").Value = -0.1652581622154127
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.002369369501396157
$ws.Range("E3").Value = 0.00008467925534031512
$ws.Range("F3").Value = 0.2844015208501288
$ws.Range("G3").Value = -0.0004875614275482367
$ws.Range("H3").Value = -0.0002308575196248227
$ws.Range("I3").Value = 0.0003402204221353358
$ws.Range("J3").Value = 0.0001211958291220369
$ws.Range("K3").Value = -0.000634069503543657
$ws.Range("L3").Value = -0.0003371146942535636
$ws.Range("M3").Value = 0.0004722282032739412
$ws.Range("N3").Value = 0.0002513753691377797
$ws.Range("O3").Value = -0.000686895551804251
$ws.Range("P3").Value = 0.0003503091363556556
$ws.Range("Q3").Value = -0.00006915549790652933
$ws.Range("R3").Value = 0.001265568880261699
$ws.Range("B4").Value = -0.002523758263964982
$ws.Range("C4").Value = 0.002369369501396157
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.1349809246210455
$ws.Range("F4").Value = 0.157923111777737
$ws.Range("G4").Value = 0.0004048714336955325
$ws.Range("H4").Value = 0.0002323913823368846
$ws.Range("I4").Value = -0.0003424809134566106
$ws.Range("J4").Value = -0.0001220010780195644
$ws.Range("K4").Value = 0.0006382823858887653
$ws.Range("L4").Value = 0.0003393545505090623
$ws.Range("M4").Value = -0.0004753657802256805
$ws.Range("N4").Value = -0.0002552659143308778
$ws.Range("O4").Value = 0.0002627615482830759
$ws.Range("P4").Value = -0.0003439280965548492
$ws.Range("Q4").Value = 0.0000239965511444554
$ws.Range("R4").Value = 0.2585231779625458
$ws.Range("B5").Value = 0.0002931071949994339
$ws.Range("C5").Value = 0.00008467925534031512
$ws.Range("D5").Value = 0.1349809246210455
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.2261890848751147
$ws.Range("G5").Value = 0.0001495056164569069
$ws.Range("H5").Value = 0.00008451679271108418
$ws.Range("I5").Value = -0.0001245544825252122
$ws.Range("J5").Value = -0.00004436971680226623
$ws.Range("K5").Value = 0.0002321324463804022
$ws.Range("L5").Value = 0.0001234174775012411
$ws.Range("M5").Value = -0.0001728824481591328
$ws.Range("N5").Value = -0.00009817768454003779
$ws.Range("O5").Value = 0.0003146591035266696
$ws.Range("P5").Value = -0.00009939581574887584
$ws.Range("Q5").Value = 0.00002389315715687617
$ws.Range("R5").Value = 0.1142621748737162
$ws.Range("B6").Value = 0.001124234086298129
$ws.Range("C6").Value = 0.2844015208501288
$ws.Range("D6").Value = 0.157923111777737
$ws.Range("E6").Value = 0.2261890848751147
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = -0.218824225336559
$ws.Range("H6").Value = 0.03486044852935592
$ws.Range("I6").Value = 0.06371089685032538
$ws.Range("J6").Value = 0.2400602062985029
$ws.Range("K6").Value = -0.2194358192293704
$ws.Range("L6").Value = 0.02395877534143304
$ws.Range("M6").Value = 0.07431340951733358
$ws.Range("N6").Value = 0.4941702061229365
$ws.Range("O6").Value = 0.08661242330160181
$ws.Range("P6").Value = -0.01860508477477728
$ws.Range("Q6").Value = 0.01058778278023932
$ws.Range("R6").Value = 0.295895750700426
$ws.Range("B7").Value = -0.00127753765741268
$ws.Range("C7").Value = -0.0004875614275482367
$ws.Range("D7").Value = 0.0004048714336955325
$ws.Range("E7").Value = 0.0001495056164569069
$ws.Range("F7").Value = -0.218824225336559
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.03486926252118899
$ws.Range("I7").Value = -0.03915802133205065
$ws.Range("J7").Value = -0.1931068753440004
$ws.Range("K7").Value = 0.9489195226014283
$ws.Range("L7").Value = -0.05549788566449686
$ws.Range("M7").Value = -0.1561456386271095
$ws.Range("N7").Value = -0.4892356764397228
$ws.Range("O7").Value = -0.1277065913529247
$ws.Range("P7").Value = -0.03738881415263039
$ws.Range("Q7").Value = -0.09807475869789299
$ws.Range("R7").Value = -0.07554794405995992
$ws.Range("B8").Value = -0.0007990859124778034
$ws.Range("C8").Value = -0.0002308575196248227
$ws.Range("D8").Value = 0.0002323913823368846
$ws.Range("E8").Value = 0.00008451679271108418
$ws.Range("F8").Value = 0.03486044852935592
$ws.Range("G8").Value = 0.03486926252118899
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0.04572593915546183
$ws.Range("J8").Value = -0.01910443812119384
$ws.Range("K8").Value = 0.05654142061270117
$ws.Range("L8").Value = -0.0457060209206458
$ws.Range("M8").Value = 0.01227609320113575
$ws.Range("N8").Value = -0.03347136095775486
$ws.Range("O8").Value = -0.03184451706527877
$ws.Range("P8").Value = 0.003501428053857957
$ws.Range("Q8").Value = -0.03663625454043101
$ws.Range("R8").Value = 0.04010910537960239
$ws.Range("B9").Value = 0.001177632623395438
$ws.Range("C9").Value = 0.0003402204221353358
$ws.Range("D9").Value = -0.0003424809134566106
$ws.Range("E9").Value = -0.0001245544825252122
$ws.Range("F9").Value = 0.06371089685032538
$ws.Range("G9").Value = -0.03915802133205065
$ws.Range("H9").Value = 0.04572593915546183
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = -0.06501116173619884
$ws.Range("K9").Value = -0.04828296343645306
$ws.Range("L9").Value = -0.03433711244891577
$ws.Range("M9").Value = 0.2648506816746283
$ws.Range("N9").Value = 0.07533673899990148
$ws.Range("O9").Value = -0.01887880047468948
$ws.Range("P9").Value = 0.2441026047930702
$ws.Range("Q9").Value = -0.001727680259869888
$ws.Range("R9").Value = 0.5525381600036388
$ws.Range("B10").Value = 0.0004195049823813824
$ws.Range("C10").Value = 0.0001211958291220369
$ws.Range("D10").Value = -0.0001220010780195644
$ws.Range("E10").Value = -0.00004436971680226623
$ws.Range("F10").Value = 0.2400602062985029
$ws.Range("G10").Value = -0.1931068753440004
$ws.Range("H10").Value = -0.01910443812119384
$ws.Range("I10").Value = -0.06501116173619884
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = -0.1809546927560668
$ws.Range("L10").Value = 0.1314273997117306
$ws.Range("M10").Value = 0.4178691396851781
$ws.Range("N10").Value = 0.2171975239587932
$ws.Range("O10").Value = 0.06950850495401474
$ws.Range("P10").Value = 0.4604969130392716
$ws.Range("Q10").Value = -0.08562035705608291
$ws.Range("R10").Value = -0.04716179437360898
$ws.Range("B11").Value = -0.00219475635291615
$ws.Range("C11").Value = -0.000634069503543657
$ws.Range("D11").Value = 0.0006382823858887653
$ws.Range("E11").Value = 0.0002321324463804022
$ws.Range("F11").Value = -0.2194358192293704
$ws.Range("G11").Value = 0.9489195226014283
$ws.Range("H11").Value = 0.05654142061270117
$ws.Range("I11").Value = -0.04828296343645306
$ws.Range("J11").Value = -0.1809546927560668
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = -0.06657911034639749
$ws.Range("M11").Value = -0.1470732051773424
$ws.Range("N11").Value = -0.4919335190286616
$ws.Range("O11").Value = -0.1436981834542915
$ws.Range("P11").Value = -0.02765637501747562
$ws.Range("Q11").Value = -0.1004173198479076
$ws.Range("R11").Value = -0.08172527082725493
$ws.Range("B12").Value = -0.001166882514834089
$ws.Range("C12").Value = -0.0003371146942535636
$ws.Range("D12").Value = 0.0003393545505090623
$ws.Range("E12").Value = 0.0001234174775012411
$ws.Range("F12").Value = 0.02395877534143304
$ws.Range("G12").Value = -0.05549788566449686
$ws.Range("H12").Value = -0.0457060209206458
$ws.Range("I12").Value = -0.03433711244891577
$ws.Range("J12").Value = 0.1314273997117306
$ws.Range("K12").Value = -0.06657911034639749
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.1708438506437887
$ws.Range("N12").Value = 0.02695412148024802
$ws.Range("O12").Value = -0.04488445878777349
$ws.Range("P12").Value = 0.09922184575571608
$ws.Range("Q12").Value = -0.03472797482349419
$ws.Range("R12").Value = 0.02011168465698987
$ws.Range("B13").Value = 0.001634561894822562
$ws.Range("C13").Value = 0.0004722282032739412
$ws.Range("D13").Value = -0.0004753657802256805
$ws.Range("E13").Value = -0.0001728824481591328
$ws.Range("F13").Value = 0.07431340951733358
$ws.Range("G13").Value = -0.1561456386271095
$ws.Range("H13").Value = 0.01227609320113575
$ws.Range("I13").Value = 0.2648506816746283
$ws.Range("J13").Value = 0.4178691396851781
$ws.Range("K13").Value = -0.1470732051773424
$ws.Range("L13").Value = 0.1708438506437887
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 0.1418231427009689
$ws.Range("O13").Value = -0.07745542996770767
$ws.Range("P13").Value = 0.4082729894982228
$ws.Range("Q13").Value = 0.008567290892541542
$ws.Range("R13").Value = 0.1355845261914255
$ws.Range("B14").Value = 0.0008368912573226979
$ws.Range("C14").Value = 0.0002513753691377797
$ws.Range("D14").Value = -0.0002552659143308778
$ws.Range("E14").Value = -0.00009817768454003779
$ws.Range("F14").Value = 0.4941702061229365
$ws.Range("G14").Value = -0.4892356764397228
$ws.Range("H14").Value = -0.03347136095775486
$ws.Range("I14").Value = 0.07533673899990148
$ws.Range("J14").Value = 0.2171975239587932
$ws.Range("K14").Value = -0.4919335190286616
$ws.Range("L14").Value = 0.02695412148024802
$ws.Range("M14").Value = 0.1418231427009689
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0.1013953598326501
$ws.Range("P14").Value = 0.08252353942066859
$ws.Range("Q14").Value = 0.08541471848978671
$ws.Range("R14").Value = 0.1447960103534467
$ws.Range("B15").Value = -0.001468032878502348
$ws.Range("C15").Value = -0.000686895551804251
$ws.Range("D15").Value = 0.0002627615482830759
$ws.Range("E15").Value = 0.0003146591035266696
$ws.Range("F15").Value = 0.08661242330160181
$ws.Range("G15").Value = -0.1277065913529247
$ws.Range("H15").Value = -0.03184451706527877
$ws.Range("I15").Value = -0.01887880047468948
$ws.Range("J15").Value = 0.06950850495401474
$ws.Range("K15").Value = -0.1436981834542915
$ws.Range("L15").Value = -0.04488445878777349
$ws.Range("M15").Value = -0.07745542996770767
$ws.Range("N15").Value = 0.1013953598326501
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = -0.1618807943276877
$ws.Range("Q15").Value = -0.00008866174028652279
$ws.Range("R15").Value = -0.000304872512271589
$ws.Range("B16").Value = 0.001248231948587792
$ws.Range("C16").Value = 0.0003503091363556556
$ws.Range("D16").Value = -0.0003439280965548492
$ws.Range("E16").Value = -0.00009939581574887584
$ws.Range("F16").Value = -0.01860508477477728
$ws.Range("G16").Value = -0.03738881415263039
$ws.Range("H16").Value = 0.003501428053857957
$ws.Range("I16").Value = 0.2441026047930702
$ws.Range("J16").Value = 0.4604969130392716
$ws.Range("K16").Value = -0.02765637501747562
$ws.Range("L16").Value = 0.09922184575571608
$ws.Range("M16").Value = 0.4082729894982228
$ws.Range("N16").Value = 0.08252353942066859
$ws.Range("O16").Value = -0.1618807943276877
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 0.0000594956240911887
$ws.Range("R16").Value = -0.01560245657341144
$ws.Range("B17").Value = -0.0001451763067674667
$ws.Range("C17").Value = -0.00006915549790652933
$ws.Range("D17").Value = 0.0000239965511444554
$ws.Range("E17").Value = 0.00002389315715687617
$ws.Range("F17").Value = 0.01058778278023932
$ws.Range("G17").Value = -0.09807475869789299
$ws.Range("H17").Value = -0.03663625454043101
$ws.Range("I17").Value = -0.001727680259869888
$ws.Range("J17").Value = -0.08562035705608291
$ws.Range("K17").Value = -0.1004173198479076
$ws.Range("L17").Value = -0.03472797482349419
$ws.Range("M17").Value = 0.008567290892541542
$ws.Range("N17").Value = 0.08541471848978671
$ws.Range("O17").Value = -0.00008866174028652279
$ws.Range("P17").Value = 0.0000594956240911887
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 0.05815244118627622
$ws.Range("B18").Value = 0.002522203803982969
$ws.Range("C18").Value = 0.001265568880261699
$ws.Range("D18").Value = 0.2585231779625458
$ws.Range("E18").Value = 0.1142621748737162
$ws.Range("F18").Value = 0.295895750700426
$ws.Range("G18").Value = -0.07554794405995992
$ws.Range("H18").Value = 0.04010910537960239
$ws.Range("I18").Value = 0.5525381600036388
$ws.Range("J18").Value = -0.04716179437360898
$ws.Range("K18").Value = -0.08172527082725493
$ws.Range("L18").Value = 0.02011168465698987
$ws.Range("M18").Value = 0.1355845261914255
$ws.Range("N18").Value = 0.1447960103534467
$ws.Range("O18").Value = -0.0003305608357407635
$ws.Range("P18").Value = -0.01560245657341144
$ws.Range("Q18").Value = 0.05815244118627622
$ws.Range("R18").Value = 1
